# Applies the "gh-pages output generated" data refresh:
#  - Sheet 1 "展览" (exhibitions) and Sheet 4 "全部类型" (all types) both
#    contain the same underlying event rows; bump the "想去人数" (interest
#    count, column F) figures that changed between scrapes, and refresh
#    the venue address + cover image for the "次元日记动漫游戏嘉年华" event.

$wb = $excel.ActiveWorkbook

$sheetExhibitions = $wb.Worksheets.Item(1)   # 展览
$sheetAllTypes    = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1: 展览 -----------------------------------------------------
$sheetExhibitions.Range("F3").Value  = 754
$sheetExhibitions.Range("F4").Value  = 1505
$sheetExhibitions.Range("F5").Value  = 230
$sheetExhibitions.Range("D6").Value  = "阜阳路16号 银瑞林国际大酒店"
$sheetExhibitions.Range("I6").Value  = "//i2.hdslb.com/bfs/openplatform/202406/LU3NYF6W1719390824853.jpeg"
$sheetExhibitions.Range("F7").Value  = 148
$sheetExhibitions.Range("F8").Value  = 6254
$sheetExhibitions.Range("F11").Value = 116
$sheetExhibitions.Range("F12").Value = 5231
$sheetExhibitions.Range("F15").Value = 1188
$sheetExhibitions.Range("F16").Value = 1188
$sheetExhibitions.Range("F17").Value = 58
$sheetExhibitions.Range("F19").Value = 69
$sheetExhibitions.Range("F21").Value = 301
$sheetExhibitions.Range("F23").Value = 3732
$sheetExhibitions.Range("F24").Value = 156

# --- Sheet 4: 全部类型 --------------------------------------------------
$sheetAllTypes.Range("F2").Value  = 84
$sheetAllTypes.Range("F4").Value  = 754
$sheetAllTypes.Range("F5").Value  = 1505
$sheetAllTypes.Range("F6").Value  = 230
$sheetAllTypes.Range("D7").Value  = "阜阳路16号 银瑞林国际大酒店"
$sheetAllTypes.Range("I7").Value  = "//i2.hdslb.com/bfs/openplatform/202406/LU3NYF6W1719390824853.jpeg"
$sheetAllTypes.Range("F8").Value  = 148
$sheetAllTypes.Range("F9").Value  = 6254
$sheetAllTypes.Range("F12").Value = 116
$sheetAllTypes.Range("F13").Value = 5231
$sheetAllTypes.Range("F16").Value = 1188
$sheetAllTypes.Range("F17").Value = 1188
$sheetAllTypes.Range("F18").Value = 58
$sheetAllTypes.Range("F20").Value = 69
$sheetAllTypes.Range("F22").Value = 301
$sheetAllTypes.Range("F24").Value = 3732
$sheetAllTypes.Range("F26").Value = 156
